# Portfolio workbook update:
#  - Positions sheet: a new broker column (IBKR) is inserted after the
#    ticker column, two new market-value columns are appended, and a new
#    GME position row is added.
#  - Trades sheet: a new GME "Buy" trade row is appended, and the Trades
#    tab becomes the active/selected sheet (Positions loses that flag).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Positions")
$ws3 = $wb.Worksheets.Item("Trades")

# --- Positions row 2 (AAPL): re-point to the new 9-column layout ---
# ticker(A) | broker(B,new) | pos_amount(C) | total_cost(D) | average_cost(E)
# | price(F,new) | market_value(G,new) | unrealized_pnl(H,new) | realized_pnl(I)
$ws1.Cells.Item(2, 1).Value = "AAPL"
$ws1.Cells.Item(2, 2).Value = "IBKR"
$ws1.Cells.Item(2, 3).Value = 20
$ws1.Cells.Item(2, 4).Value = 2611
$ws1.Cells.Item(2, 5).Value = 130.55
$ws1.Cells.Item(2, 6).Value = 190.6900024414062
$ws1.Cells.Item(2, 7).Value = 3813.800048828125
$ws1.Cells.Item(2, 8).Value = 1202.800048828125
$ws1.Cells.Item(2, 9).Value = 8367

# --- Positions row 3 (new GME position) ---
$ws1.Cells.Item(3, 1).Value = "GME"
$ws1.Cells.Item(3, 2).Value = "IBKR"
$ws1.Cells.Item(3, 3).Value = 10
$ws1.Cells.Item(3, 4).Value = 1002
$ws1.Cells.Item(3, 5).Value = 0
$ws1.Cells.Item(3, 6).Value = 22.96999931335449
$ws1.Cells.Item(3, 7).Value = 229.6999931335449
$ws1.Cells.Item(3, 8).Value = -772.3000068664551
$ws1.Cells.Item(3, 9).Value = 0

# --- Trades: append the GME buy trade on row 5 ---
$ws3.Cells.Item(4, 1).Copy($ws3.Cells.Item(5, 1)) | Out-Null   # reuse row 4's date style (m/d/yyyy)
$ws3.Cells.Item(5, 1).Value = 45124
$ws3.Cells.Item(5, 2).Value = "GameStop"
$ws3.Cells.Item(5, 3).Value = "GME"
$ws3.Cells.Item(5, 4).Value = "Buy"
$ws3.Cells.Item(5, 5).Value = 10
$ws3.Cells.Item(5, 6).Value = 100
$ws3.Cells.Item(5, 7).Formula = "=F5*E5"
$ws3.Cells.Item(5, 8).Value = 2
$ws3.Cells.Item(5, 9).Value = "Speculative"

# --- Active sheet / selection bookkeeping ---
# Trades becomes the active tab with H6 selected; Positions keeps C6 selected
# but is no longer the active tab.
$ws3.Activate() | Out-Null
$ws3.Range("H6").Select() | Out-Null

Write-Output "portfolio updated"
